$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.622.39'
$ws.Range("E2").Value = '  +1.11%  '
$ws.Range("D3").Value = '1.558.57'
$ws.Range("E3").Value = '  -0.98%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.38'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  -0.71%  '
$ws.Range("E6").Value = '  -0.72%  '
$ws.Range("E7").Value = '  -0.29%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.60'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  +3.55%  '
$ws.Range("E9").Value = '  -0.22%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0893'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  -0.13%  '
$ws.Range("D12").Value = '1.779.85'
$ws.Range("E12").Value = '  -1.10%  '
$ws.Range("D13").Value = '1.556.71'
$ws.Range("E13").Value = '  -0.77%  '
$ws.Range("D14").Value = '28.626.56'
$ws.Range("E14").Value = '  +1.04%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.515'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  -0.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.63'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  -1.32%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.43'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  -0.25%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.89'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.37'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  -0.60%  '
$ws.Range("E20").Value = '  -1.89%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.998'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  -0.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.91'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  -0.96%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.98'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  -0.66%  '
$ws.Range("E24").Value = '  +1.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.06'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  -0.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14.76'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  -1.20%  '
$ws.Range("E27").Value = '  -0.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  -0.26%  '
$ws.Range("E29").Value = '  -2.02%  '
$ws.Range("E30").Value = '  -4.06%  '
$ws.Range("E31").Value = '  -1.78%  '
$ws.Range("E32").Value = '  -0.87%  '
$ws.Range("D33").Value = '1.391.15'
$ws.Range("E34").Value = '  -2.64%  '
$ws.Range("E35").Value = '  -2.82%  '
$ws.Range("E36").Value = '  -1.83%  '
$ws.Range("E37").Value = '  +0.56%  '
$ws.Range("E38").Value = '  -3.70%  '
$ws.Range("E39").Value = '  -0.70%  '
$ws.Range("E40").Value = '  +3.49%  '
$ws.Range("E41").Value = '  -0.42%  '
$ws.Range("E42").Value = '  -0.28%  '
$ws.Range("E43").Value = '  -1.23%  '
$ws.Range("E44").Value = '  +0.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.96'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  +2.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.28'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  -1.86%  '
$ws.Range("D47").Value = '1.693.32'
$ws.Range("E47").Value = '  -1.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.869'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  -5.99%  '
$ws.Range("E49").Value = '  -0.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '43.31'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +4.75%  '
$ws.Range("E51").Value = '  -0.46%  '
